# Rename the "_old"/"_new" column-header suffixes to the respective
# input-file-name suffixes ("_FV2310" / "_FV2404") and turn the sheet's
# data range into a proper Excel Table, with the header row frozen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) -----------------------------------
$lastCol = 21   # A..U
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old = $cell.Value2
    if ($old -ne $null) {
        if ($old -like "*_old") {
            $cell.Value2 = ($old -replace "_old$", "_FV2310")
        } elseif ($old -like "*_new") {
            $cell.Value2 = ($old -replace "_new$", "_FV2404")
        }
    }
}

# --- 2. Freeze the header row (pane split after row 1) -----------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table (ListObject) -----------
$rng = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
